$d = $word.ActiveDocument

# 1. Repoint the first hyperlink ("Google Java Style Guide") from rId12
#    (Google Java style guide) to rId13 (swagger.io) - reuses the existing
#    swagger.io external relationship, same as the target diff.
$h = $d.Hyperlinks.Item(1)
$h.Address = "http://swagger.io/specification/"

# 2. The empty paragraph right after that hyperlink paragraph gains the
#    "_GoBack" bookmark (it moves here from its old home inside the
#    REST API section, which is being removed below).
$p7 = $d.Paragraphs.Item(7)
$d.Bookmarks.Add("_GoBack", $p7.Range)

# 3. Revert the whole "REST API" naming-standard section: delete every
#    paragraph from the "REST API" Heading1 through the end of the
#    document body (the content added by the reverted commit).
$first = $d.Paragraphs.Item(14)
$deleteRange = $d.Range($first.Range.Start, $d.Content.End)
$deleteRange.Delete()
